# Scheme of Work document update
# - Table 1 (info table): add Province/District/Sector/RQF Level rows,
#   update School/Module/Trainer text, remove Term row
# - Table 2 (scheme grid): update first data row, add two more term rows

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: set a table cell's text, turning a "|||" marker into a real
# Word line break (<w:br/>) using a unique Find/Replace so it can't
# collide with text elsewhere in the document.
# ---------------------------------------------------------------------
function Set-CellText($table, $row, $col, $text) {
    $parts = $text -split "\|\|\|"
    $table.Cell($row, $col).Range.Text = $parts -join "<<WORDBREAK>>"
    if ($parts.Length -gt 1) {
        $cellRange = $table.Cell($row, $col).Range
        $cellRange.Find.Execute("<<WORDBREAK>>", $false, $false, $false, $false, $false, $true, 1, $false, "^l", 2) | Out-Null
    }
}

# =======================================================================
# Table 1: header / info table
# =======================================================================
$t1 = $d.Tables.Item(1)

# Current rows: 1 School:, 2 Module Code & Title:, 3 Trainer Name:,
#               4 School Year:, 5 Term:

# Insert "Province:" / "District:" rows before the current row 1 (School:)
$t1.Rows.Add($t1.Rows.Item(1)) | Out-Null
$t1.Rows.Add($t1.Rows.Item(1)) | Out-Null
Set-CellText $t1 1 1 "Province:"
Set-CellText $t1 1 2 "Kigali City"
Set-CellText $t1 2 1 "District:"
Set-CellText $t1 2 2 "Gasabo"

# Row 3 is now "School:" -> update value, and insert "Sector:" row after it
Set-CellText $t1 3 2 "IPRC Kigali"
$t1.Rows.Add($t1.Rows.Item(4)) | Out-Null
Set-CellText $t1 4 1 "Sector:"
Set-CellText $t1 4 2 "ICT & MULTIMEDIA"

# Row 5 is now "Module Code & Title:" -> update value, insert "RQF Level:" after it
Set-CellText $t1 5 2 "CSAPA301: C Programming"
$t1.Rows.Add($t1.Rows.Item(6)) | Out-Null
Set-CellText $t1 6 1 "RQF Level:"
Set-CellText $t1 6 2 "Level 4"

# Row 7 is now "Trainer Name:" -> update value
Set-CellText $t1 7 2 "TUYISINGIZE Leonard"

# Row 8 is "School Year:" (unchanged), Row 9 is "Term:" -> delete it
$t1.Rows.Item(9).Delete()

# =======================================================================
# Table 2: scheme-of-work grid
# =======================================================================
$t2 = $d.Tables.Item(2)

# Row 3 holds the first (and, originally, only) data row; update in place.
Set-CellText $t2 3 1 "Week 1-12 (Jan-Mar 2025)"
Set-CellText $t2 3 2 "LO1: Understand C basics|||LO2: Write simple programs"
Set-CellText $t2 3 3 "40 hours"
Set-CellText $t2 3 4 "IC1.1: Variables and data types|||IC1.2: Control structures"
Set-CellText $t2 3 5 "Blended learning"
Set-CellText $t2 3 6 "Computers, C compiler, textbooks"
Set-CellText $t2 3 7 "Weekly quizzes, practical exercises"
Set-CellText $t2 3 8 "Lab/Classroom"
Set-CellText $t2 3 9 "Term 1"

# New row: Week 13-24
$t2.Rows.Add() | Out-Null
$row = $t2.Rows.Count
Set-CellText $t2 $row 1 "Week 13-24 (Apr-Jun 2025)"
Set-CellText $t2 $row 2 "LO3: Work with arrays|||LO4: Use functions"
Set-CellText $t2 $row 3 "40 hours"
Set-CellText $t2 $row 4 "IC2.1: Arrays and strings|||IC2.2: Functions and scope"
Set-CellText $t2 $row 5 "Blended learning"
Set-CellText $t2 $row 6 "Computers, C compiler, textbooks"
Set-CellText $t2 $row 7 "Weekly quizzes, practical exercises"
Set-CellText $t2 $row 8 "Lab/Classroom"
Set-CellText $t2 $row 9 "Term 2"

# New row: Week 25-36
$t2.Rows.Add() | Out-Null
$row = $t2.Rows.Count
Set-CellText $t2 $row 1 "Week 25-36 (Jul-Sep 2025)"
Set-CellText $t2 $row 2 "LO5: Manage pointers|||LO6: Handle files"
Set-CellText $t2 $row 3 "40 hours"
Set-CellText $t2 $row 4 "IC3.1: Pointers and memory|||IC3.2: File operations"
Set-CellText $t2 $row 5 "Blended learning"
Set-CellText $t2 $row 6 "Computers, C compiler, textbooks"
Set-CellText $t2 $row 7 "Weekly quizzes, practical exercises"
Set-CellText $t2 $row 8 "Lab/Classroom"
Set-CellText $t2 $row 9 "Term 3"

Write-Host "Done."
